$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a flat positional dump (Conta / Nome / Saldo) with no merged
# cells or formulas. The update adds several new account rows and removes
# one stale row. Operations are applied from the bottom of the sheet upward
# so that row numbers referenced below (taken from the original layout)
# remain valid as each step executes.
# ---------------------------------------------------------------------------

# 1) Remove the old "004452597 / LARA / 36.17" row (row 97 in the original
#    layout). Its account now reappears higher up the sheet with a new
#    balance, so the old low-balance entry is deleted outright.
$ws.Rows.Item(97).Delete()

# 2) Insert "005440756 / VALERIA / 998.57" right before the
#    "001651617 / MIRELLA / 931.64" row (row 10).
$ws.Rows.Item(10).Insert()
$ws.Cells.Item(10, 1).Value = "'005440756"
$ws.Cells.Item(10, 2).Value = "VALERIA"
$ws.Cells.Item(10, 3).Value = 998.57

# 3) Insert "004267119 / ANA / 1494.47" right before the
#    "005143579 / GABRIEL / 1369.2" row (row 9).
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "'004267119"
$ws.Cells.Item(9, 2).Value = "ANA"
$ws.Cells.Item(9, 3).Value = 1494.47

# 4) Insert "004461526 / ASSISTIGAS / 4984.76" right before the
#    "004313254 / GUSTAVO / 4292" row (row 7).
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value = "'004461526"
$ws.Cells.Item(7, 2).Value = "ASSISTIGAS"
$ws.Cells.Item(7, 3).Value = 4984.76

# 5) Row 5 used to hold "004693349 / CATARINE / 12691.6"; that account was
#    replaced with a different one, so overwrite the row in place.
#    (A leading apostrophe forces the account number to stay text, keeping
#    the leading zeros intact instead of being parsed as a number.)
$ws.Cells.Item(5, 1).Value = "'004398253"
$ws.Cells.Item(5, 2).Value = "EULER"
$ws.Cells.Item(5, 3).Value = 19127.11

# 6) Insert "004940699 / RACHEL / 18534.56" right after the row above.
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).Value = "'004940699"
$ws.Cells.Item(6, 2).Value = "RACHEL"
$ws.Cells.Item(6, 3).Value = 18534.56

# 7) Insert "004452597 / LARA / 13175.78" right after that (this is the same
#    account number that was removed from further down the sheet, now with
#    an updated balance).
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value = "'004452597"
$ws.Cells.Item(7, 2).Value = "LARA"
$ws.Cells.Item(7, 3).Value = 13175.78
